$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit swaps the two observation records currently sitting in rows 4
# and 5 of the "Artfynd" sheet (same species-occurrence table, two rows
# trade places). Columns C, P, S, T, U, V, W, AD, AE, AG, AI, AT, AW, AX,
# AY already hold identical values in both rows, so they are left as-is.
#
# Reads use .Value2 (this host's .Value getter doesn't collapse to a plain
# scalar when captured into a variable). Writes to the "Antal" (I),
# "Startdatum" (Y) and "Slutdatum" (AA) cells temporarily set NumberFormat
# to "@" (Text) before assigning, because those values look numeric/
# date-like and would otherwise be auto-parsed into a number/date serial
# instead of staying literal text (same as typing them straight into
# Excel would do). The style is put back to "Normal" right after so the
# cell ends up with its original default formatting, only the text value
# having changed.

# --- Capture current (pre-edit) row 4 values ---
$A4 = $ws.Range("A4").Value2
$B4 = $ws.Range("B4").Value2
$D4 = $ws.Range("D4").Value2
$E4 = $ws.Range("E4").Value2
$F4 = $ws.Range("F4").Value2
$G4 = $ws.Range("G4").Value2
$H4 = $ws.Range("H4").Value2
$I4 = $ws.Range("I4").Value2
$J4 = $ws.Range("J4").Value2
$K4 = $ws.Range("K4").Value2
$N4 = $ws.Range("N4").Value2
$Q4 = $ws.Range("Q4").Value2
$R4 = $ws.Range("R4").Value2
$Y4 = $ws.Range("Y4").Value2
$AA4 = $ws.Range("AA4").Value2
$AO4 = $ws.Range("AO4").Value2
$AQ4 = $ws.Range("AQ4").Value2
$AR4 = $ws.Range("AR4").Value2

# --- Capture current (pre-edit) row 5 values ---
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$D5 = $ws.Range("D5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$M5 = $ws.Range("M5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2
$Y5 = $ws.Range("Y5").Value2
$AA5 = $ws.Range("AA5").Value2
$AO5 = $ws.Range("AO5").Value2

# --- Row 4 becomes what row 5 used to hold ---
$ws.Range("A4").Value = $A5
$ws.Range("B4").Value = $B5
$ws.Range("D4").Value = $D5
$ws.Range("E4").Value = $E5
$ws.Range("F4").Value = $F5
$ws.Range("G4").Value = $G5
$ws.Range("H4").Value = $H5
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("M4").Value = $M5
$ws.Range("Q4").Value = $Q5
$ws.Range("R4").Value = $R5
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = $Y5
$ws.Range("Y4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = $AA5
$ws.Range("AA4").Style = "Normal"
$ws.Range("AO4").Value = $AO5
$ws.Range("AQ4").Value = ""
$ws.Range("AR4").Value = ""

# --- Row 5 becomes what row 4 used to hold ---
$ws.Range("A5").Value = $A4
$ws.Range("B5").Value = $B4
$ws.Range("D5").Value = $D4
$ws.Range("E5").Value = $E4
$ws.Range("F5").Value = $F4
$ws.Range("G5").Value = $G4
$ws.Range("H5").Value = $H4
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = $I4
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = $J4
$ws.Range("K5").Value = $K4
$ws.Range("N5").Value = $N4
$ws.Range("M5").Value = ""
$ws.Range("Q5").Value = $Q4
$ws.Range("R5").Value = $R4
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = $Y4
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = $AA4
$ws.Range("AA5").Style = "Normal"
$ws.Range("AO5").Value = $AO4
$ws.Range("AQ5").Value = $AQ4
$ws.Range("AR5").Value = $AR4
